$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two input values; dependent formulas (F4, F5) recalc automatically.
$ws.Range("F2").Value = 929705
$ws.Range("F3").Value = 912245

# Update selected cell to match the saved cursor position.
$ws.Range("F9").Select()
